$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from the default "Sheet1" to something that
# reflects what it actually holds: chinook sales data.
$ws.Name = "sales_data"
